$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 4 (old row4 -> row5, old row5 -> row6)
# ---------------------------------------------------------------------------
$ws.Rows("4").Insert()

# ---------------------------------------------------------------------------
# 2. Update text content
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in C. Poti Municipality"
$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("A5").Value = "disabilities Persons "

# ---------------------------------------------------------------------------
# 3. Update numeric data
# ---------------------------------------------------------------------------
$cols = @("B","C","D","E","F","G","H","I")

$row4vals = @(436,411,376,379,371,370,359,365)
for ($i=0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "4").Value = $row4vals[$i]
}

$row5vals = @(466,442,405,404,395,396,381,389)
for ($i=0; $i -lt $cols.Length; $i++) {
  $ws.Range($cols[$i] + "5").Value = $row5vals[$i]
}

# ---------------------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------------------
$ws.Rows("1").RowHeight = 51
$ws.Rows("4").RowHeight = 24.75
$ws.Rows("5").RowHeight = 21
$ws.Rows("6").RowHeight = 27.75

# ---------------------------------------------------------------------------
# 5. Column A width
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20

# ---------------------------------------------------------------------------
# 6. Merge title row A1:I1
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").Merge()

# ---------------------------------------------------------------------------
# 7. Formatting - Title row (A1:I1)
# ---------------------------------------------------------------------------
$ws.Range("A1:I1").HorizontalAlignment = -4108
$ws.Range("A1:I1").VerticalAlignment = -4108
$ws.Range("A1:I1").WrapText = $true

# ---------------------------------------------------------------------------
# 8. Formatting - A3 (blank corner cell, changes to Sylfaen 11)
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Name = "Sylfaen"
$ws.Range("A3").Font.Size = 11
$ws.Range("A3").Font.ThemeColor = 1

# ---------------------------------------------------------------------------
# 9. Formatting - A4 (family with disabilities Persons)
# ---------------------------------------------------------------------------
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").Font.ThemeColor = 1
$ws.Range("A4").Interior.Pattern = 1
$ws.Range("A4").Interior.ThemeColor = 0
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("A4").Borders.Item(8).Weight = 2
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4108
$ws.Range("A4").WrapText = $true

# ---------------------------------------------------------------------------
# 10. Formatting - A5 (disabilities Persons)
# ---------------------------------------------------------------------------
$ws.Range("A5").Font.Name = "Arial"
$ws.Range("A5").Font.Size = 10
$ws.Range("A5").Font.ThemeColor = 1
$ws.Range("A5").Borders.Item(8).LineStyle = 0
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("A5").Borders.Item(9).Weight = 2
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A5").WrapText = $true

# ---------------------------------------------------------------------------
# 11. Formatting - A6 (source row)
# ---------------------------------------------------------------------------
$ws.Range("A6").Borders.Item(8).LineStyle = 0
$ws.Range("A6").WrapText = $true
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("A6").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 12. Formatting - value rows (B4:I4, B5:H5) - remove right-alignment, set number format
# ---------------------------------------------------------------------------
$ws.Range("B4:I5").NumberFormat = "#\ ##0"
$ws.Range("B4:I4").HorizontalAlignment = -4142
$ws.Range("B5:I5").HorizontalAlignment = -4142

# ---------------------------------------------------------------------------
# 13. Formatting - I5 gets a bottom border only
# ---------------------------------------------------------------------------
$ws.Range("I5").Borders.Item(8).LineStyle = 0
$ws.Range("I5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").Borders.Item(9).Weight = 2
